{"js": "// Update the 100 arithmetic-problem cells in the single 5x20 table.\n// The diff only changes the <w:t> text content of each cell's single run,\n// so we replace the cell body's range text (which preserves the existing\n// run/paragraph formatting) rather than clearing + re-inserting.\n\nconst beforeValues = [\"3+51=\", \"90-68=\", \"24+49=\", \"0+62=\", \"54-51=\", \"56+22=\", \"75-11=\", \"93-70=\", \"77+5=\", \"33+36=\", \"36+24=\", \"98-43=\", \"91+8=\", \"86-1=\", \"91-58=\", \"15+83=\", \"41-16=\", \"58+20=\", \"25-15=\", \"40+14=\", \"35-15=\", \"71-59=\", \"17+16=\", \"47-7=\", \"61-50=\", \"59-8=\", \"20+4=\", \"16-3=\", \"33+14=\", \"91-80=\", \"59-56=\", \"48+2=\", \"82-50=\", \"56-19=\", \"12+16=\", \"62+31=\", \"60-59=\", \"54+35=\", \"81-7=\", \"64-19=\", \"47+23=\", \"74-56=\", \"17+19=\", \"85-85=\", \"30-15=\", \"93-83=\", \"7+29=\", \"27+25=\", \"5+59=\", \"28+7=\", \"54+29=\", \"71-37=\", \"78-10=\", \"77+20=\", \"54-23=\", \"12+74=\", \"60+25=\", \"68-12=\", \"25-9=\", \"95-15=\", \"13-6=\", \"89-69=\", \"45+9=\", \"15+31=\", \"6+9=\", \"8+87=\", \"86-63=\", \"97-1=\", \"22+58=\", \"96-58=\", \"73-20=\", \"73-69=\", \"7+32=\", \"83-30=\", \"18+37=\", \"0+19=\", \"0+25=\", \"79+15=\", \"28+62=\", \"51+14=\", \"29+63=\", \"63-47=\", \"81+13=\", \"12+39=\", \"18+32=\", \"75+13=\", \"20+6=\", \"17+55=\", \"15+25=\", \"8+15=\", \"52-46=\", \"19+22=\", \"60-22=\", \"61-19=\", \"66-52=\", \"16-0=\", \"59-42=\", \"63-55=\", \"52+13=\", \"99-49=\"];\nconst afterValues = [\"98-80=\", \"56-18=\", \"62+14=\", \"41+45=\", \"78-67=\", \"99-30=\", \"33+2=\", \"86-77=\", \"93-49=\", \"72-59=\", \"75-5=\", \"25+60=\", \"35+4=\", \"33+8=\", \"25+21=\", \"71-11=\", \"37+41=\", \"61-31=\", \"65-6=\", \"63+22=\", \"14+59=\", \"80-76=\", \"19+51=\", \"9+17=\", \"16+3=\", \"59+8=\", \"90-84=\", \"60-3=\", \"22+52=\", \"24+21=\", \"27+48=\", \"17+40=\", \"87-25=\", \"62-38=\", \"11+5=\", \"49-38=\", \"47+8=\", \"58-34=\", \"29+4=\", \"2+87=\", \"92-36=\", \"14+6=\", \"90-81=\", \"62-19=\", \"71+15=\", \"50+46=\", \"22+28=\", \"78-48=\", \"16+82=\", \"35-23=\", \"28+43=\", \"14+39=\", \"42+44=\", \"49+1=\", \"94-56=\", \"80-50=\", \"36-21=\", \"41-9=\", \"3+38=\", \"44+6=\", \"63-44=\", \"72-56=\", \"20+22=\", \"77-39=\", \"84-23=\", \"95-15=\", \"71-18=\", \"30+12=\", \"25+53=\", \"91-76=\", \"19+35=\", \"9+18=\", \"39+53=\", \"72-17=\", \"62+27=\", \"10+32=\", \"94-22=\", \"27-11=\", \"58-20=\", \"64-28=\", \"24+52=\", \"59-9=\", \"77-60=\", \"86-34=\", \"83-12=\", \"11+56=\", \"70-63=\", \"62-25=\", \"0+82=\", \"37-5=\", \"17+0=\", \"82-37=\", \"82-68=\", \"86-60=\", \"0+72=\", \"0+32=\", \"59+36=\", \"41+51=\", \"12+11=\", \"98-40=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columns = 5;\nconst rows = afterValues.length / columns;\n\nif (table.rowCount < rows) {\n  throw new Error(\"Table has fewer rows than expected: \" + table.rowCount);\n}\n\n// Load the current text of every cell so we can confirm we're editing the\n// expected \"before\" problem before overwriting it with the \"after\" one.\nconst cells = [];\nfor (let i = 0; i < afterValues.length; i++) {\n  const r = Math.floor(i / columns);\n  const c = i % columns;\n  const cell = table.getCell(r, c);\n  cell.body.load(\"text\");\n  cells.push(cell);\n}\nawait context.sync();\n\nfor (let i = 0; i < afterValues.length; i++) {\n  const cell = cells[i];\n  const current = cell.body.text.trim();\n  if (current !== beforeValues[i]) {\n    throw new Error(\n      `Unexpected text in cell ${i} (expected \"${beforeValues[i]}\", found \"${current}\")`\n    );\n  }\n  const range = cell.body.getRange();\n  range.insertText(afterValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the 100 arithmetic-problem cells in the single 5x20 table.\n# The diff only changes the text content inside each cell's single run,\n# so we assign directly to Cell.Range.Text, which replaces the run's text\n# in place and preserves the existing run/paragraph formatting (font,\n# size, alignment) rather than clearing and re-creating the cell content.\n\n$beforeValues = @('3+51=', '90-68=', '24+49=', '0+62=', '54-51=', '56+22=', '75-11=', '93-70=', '77+5=', '33+36=', '36+24=', '98-43=', '91+8=', '86-1=', '91-58=', '15+83=', '41-16=', '58+20=', '25-15=', '40+14=', '35-15=', '71-59=', '17+16=', '47-7=', '61-50=', '59-8=', '20+4=', '16-3=', '33+14=', '91-80=', '59-56=', '48+2=', '82-50=', '56-19=', '12+16=', '62+31=', '60-59=', '54+35=', '81-7=', '64-19=', '47+23=', '74-56=', '17+19=', '85-85=', '30-15=', '93-83=', '7+29=', '27+25=', '5+59=', '28+7=', '54+29=', '71-37=', '78-10=', '77+20=', '54-23=', '12+74=', '60+25=', '68-12=', '25-9=', '95-15=', '13-6=', '89-69=', '45+9=', '15+31=', '6+9=', '8+87=', '86-63=', '97-1=', '22+58=', '96-58=', '73-20=', '73-69=', '7+32=', '83-30=', '18+37=', '0+19=', '0+25=', '79+15=', '28+62=', '51+14=', '29+63=', '63-47=', '81+13=', '12+39=', '18+32=', '75+13=', '20+6=', '17+55=', '15+25=', '8+15=', '52-46=', '19+22=', '60-22=', '61-19=', '66-52=', '16-0=', '59-42=', '63-55=', '52+13=', '99-49=')\n$afterValues = @('98-80=', '56-18=', '62+14=', '41+45=', '78-67=', '99-30=', '33+2=', '86-77=', '93-49=', '72-59=', '75-5=', '25+60=', '35+4=', '33+8=', '25+21=', '71-11=', '37+41=', '61-31=', '65-6=', '63+22=', '14+59=', '80-76=', '19+51=', '9+17=', '16+3=', '59+8=', '90-84=', '60-3=', '22+52=', '24+21=', '27+48=', '17+40=', '87-25=', '62-38=', '11+5=', '49-38=', '47+8=', '58-34=', '29+4=', '2+87=', '92-36=', '14+6=', '90-81=', '62-19=', '71+15=', '50+46=', '22+28=', '78-48=', '16+82=', '35-23=', '28+43=', '14+39=', '42+44=', '49+1=', '94-56=', '80-50=', '36-21=', '41-9=', '3+38=', '44+6=', '63-44=', '72-56=', '20+22=', '77-39=', '84-23=', '95-15=', '71-18=', '30+12=', '25+53=', '91-76=', '19+35=', '9+18=', '39+53=', '72-17=', '62+27=', '10+32=', '94-22=', '27-11=', '58-20=', '64-28=', '24+52=', '59-9=', '77-60=', '86-34=', '83-12=', '11+56=', '70-63=', '62-25=', '0+82=', '37-5=', '17+0=', '82-37=', '82-68=', '86-60=', '0+72=', '0+32=', '59+36=', '41+51=', '12+11=', '98-40=')\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$columns = 5\n$rows = $afterValues.Length / $columns\n\nfor ($i = 0; $i -lt $afterValues.Length; $i++) {\n    $r = [math]::Floor($i / $columns) + 1\n    $c = ($i % $columns) + 1\n    $cell = $t.Cell($r, $c)\n\n    # Cell.Range.Text carries trailing cell-mark characters (CR + BEL);\n    # strip them before comparing against the expected \"before\" value.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $beforeValues[$i]) {\n        throw \"Unexpected text in cell $i (expected '$($beforeValues[$i])', found '$current')\"\n    }\n\n    $cell.Range.Text = $afterValues[$i]\n}\n"}
